$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a plain-text value to a cell, forcing text storage so
# numeric-looking strings (NIK/HP with leading zeros, 17+ digit IDs) are
# not silently coerced into numbers, then drop back to the default style
# so no stray number format is left on the cell.
function Set-PlainTextCell($cell, [string]$text) {
    $cell.NumberFormat = "@"
    if ($text -eq "") {
        # A bare value of empty string clears the cell instead of storing
        # an empty text value, so use a leading quote (text-qualifier)
        # to force an empty *text* cell, matching the source sheet.
        $cell.Value = ""
    } else {
        $cell.Value = $text
    }
    $cell.Style = "Normal"
}

# Helper: write the NIK column, which keeps the sheet's existing text
# number format (style index 9 in the original workbook) instead of
# reverting to the default style.
function Set-NikCell($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# Row 101: S097
Set-PlainTextCell $ws.Cells.Item(101, 1) 'S097'
Set-PlainTextCell $ws.Cells.Item(101, 2) 'S'
Set-NikCell $ws.Cells.Item(101, 3) '3309120108049001'
Set-PlainTextCell $ws.Cells.Item(101, 4) 'ALE SANDRO DWI SAPUTRA'
Set-PlainTextCell $ws.Cells.Item(101, 5) '0895340453050'
Set-PlainTextCell $ws.Cells.Item(101, 6) 'KLIWONAN RT 02/07 DESA JERON, KEC NOGOSARI, KAB BOYOLALI'
$ws.Cells.Item(101, 7).Value = 0
Set-PlainTextCell $ws.Cells.Item(101, 8) ''
Set-PlainTextCell $ws.Cells.Item(101, 9) 'static/uploads/S097.png'

# Row 102: S098
Set-PlainTextCell $ws.Cells.Item(102, 1) 'S098'
Set-PlainTextCell $ws.Cells.Item(102, 2) 'S'
Set-NikCell $ws.Cells.Item(102, 3) '3314011507800004'
Set-PlainTextCell $ws.Cells.Item(102, 4) 'LASNO'
Set-PlainTextCell $ws.Cells.Item(102, 5) '082220833940'
Set-PlainTextCell $ws.Cells.Item(102, 6) 'KEDEN WETAN RT 14 DESA KEDEN, KEC KALIJAMBE, KAB SRAGEN'
$ws.Cells.Item(102, 7).Value = 0
Set-PlainTextCell $ws.Cells.Item(102, 8) ''
Set-PlainTextCell $ws.Cells.Item(102, 9) 'static/uploads/S098.png'

# Row 103: S099
Set-PlainTextCell $ws.Cells.Item(103, 1) 'S099'
Set-PlainTextCell $ws.Cells.Item(103, 2) 'S'
Set-NikCell $ws.Cells.Item(103, 3) '3319060511840002'
Set-PlainTextCell $ws.Cells.Item(103, 4) 'MATTORI'
Set-PlainTextCell $ws.Cells.Item(103, 5) '082322083412'
Set-PlainTextCell $ws.Cells.Item(103, 6) 'RT 01/07 DESA TERBAN, KEC JEKULO, KAB KUDUS'
$ws.Cells.Item(103, 7).Value = 0
Set-PlainTextCell $ws.Cells.Item(103, 8) ''
Set-PlainTextCell $ws.Cells.Item(103, 9) 'static/uploads/S099.png'

# Row 104: S100
Set-PlainTextCell $ws.Cells.Item(104, 1) 'S100'
Set-PlainTextCell $ws.Cells.Item(104, 2) 'S'
Set-NikCell $ws.Cells.Item(104, 3) '3319080808890003'
Set-PlainTextCell $ws.Cells.Item(104, 4) 'AZWAR ANAS'
Set-PlainTextCell $ws.Cells.Item(104, 5) '085225324480'
Set-PlainTextCell $ws.Cells.Item(104, 6) 'NGARINGAN RT 06/06 DESA KLUMPIT, KEC GEBOG, KAB KUDUS'
$ws.Cells.Item(104, 7).Value = 0
Set-PlainTextCell $ws.Cells.Item(104, 8) ''
Set-PlainTextCell $ws.Cells.Item(104, 9) 'static/uploads/S100.png'

# Row 105: S101
Set-PlainTextCell $ws.Cells.Item(105, 1) 'S101'
Set-PlainTextCell $ws.Cells.Item(105, 2) 'S'
Set-NikCell $ws.Cells.Item(105, 3) '3318081006740001'
Set-PlainTextCell $ws.Cells.Item(105, 4) 'WAWAN SAFUAN'
Set-PlainTextCell $ws.Cells.Item(105, 5) '085281656966'
Set-PlainTextCell $ws.Cells.Item(105, 6) 'LANGGEN RT 01/03 DESA LANGGENHARJO, KEC JUWANA, KAB PATI'
$ws.Cells.Item(105, 7).Value = 0
Set-PlainTextCell $ws.Cells.Item(105, 8) ''
Set-PlainTextCell $ws.Cells.Item(105, 9) 'static/uploads/S101.png'

# Row 106: S102
Set-PlainTextCell $ws.Cells.Item(106, 1) 'S102'
Set-PlainTextCell $ws.Cells.Item(106, 2) 'S'
Set-NikCell $ws.Cells.Item(106, 3) '3319062305970004'
Set-PlainTextCell $ws.Cells.Item(106, 4) 'NURUL FARIHIN'
Set-PlainTextCell $ws.Cells.Item(106, 5) '082325160701'
Set-PlainTextCell $ws.Cells.Item(106, 6) 'RT 05/07 DESA BULUNG KULON, KEC JEKULO, KAB KUDUS'
$ws.Cells.Item(106, 7).Value = 0
Set-PlainTextCell $ws.Cells.Item(106, 8) ''
Set-PlainTextCell $ws.Cells.Item(106, 9) 'static/uploads/S102.png'

# Row 107: S103
Set-PlainTextCell $ws.Cells.Item(107, 1) 'S103'
Set-PlainTextCell $ws.Cells.Item(107, 2) 'S'
Set-NikCell $ws.Cells.Item(107, 3) '3173061807990008'
Set-PlainTextCell $ws.Cells.Item(107, 4) 'HARYONO'
Set-PlainTextCell $ws.Cells.Item(107, 5) '0'
Set-PlainTextCell $ws.Cells.Item(107, 6) 'RT 04/01 DESA PADANG JAYA, KEC MAJENANG, KAB CILACAP'
$ws.Cells.Item(107, 7).Value = 0
Set-PlainTextCell $ws.Cells.Item(107, 8) ''
Set-PlainTextCell $ws.Cells.Item(107, 9) 'static/uploads/S103.png'

# Row 108: S104
Set-PlainTextCell $ws.Cells.Item(108, 1) 'S104'
Set-PlainTextCell $ws.Cells.Item(108, 2) 'S'
Set-NikCell $ws.Cells.Item(108, 3) '33181219068200011'
Set-PlainTextCell $ws.Cells.Item(108, 4) 'KIYANTO'
Set-PlainTextCell $ws.Cells.Item(108, 5) '085327340315'
Set-PlainTextCell $ws.Cells.Item(108, 6) 'MAWAR RT 02/05 DESA JAMBEAN KIDUL, KEC MARGOREJO, KAB PATI'
$ws.Cells.Item(108, 7).Value = 0
Set-PlainTextCell $ws.Cells.Item(108, 8) ''
Set-PlainTextCell $ws.Cells.Item(108, 9) 'static/uploads/S104.png'

# Row 109: S105
Set-PlainTextCell $ws.Cells.Item(109, 1) 'S105'
Set-PlainTextCell $ws.Cells.Item(109, 2) 'S'
Set-NikCell $ws.Cells.Item(109, 3) '3318042004630001'
Set-PlainTextCell $ws.Cells.Item(109, 4) 'TRIYONO'
Set-PlainTextCell $ws.Cells.Item(109, 5) '082146354406'
Set-PlainTextCell $ws.Cells.Item(109, 6) 'BLLIBAK RT 02/02 DESA PULOREJO, KEC WINONG, KAB PATI'
$ws.Cells.Item(109, 7).Value = 0
Set-PlainTextCell $ws.Cells.Item(109, 8) ''
Set-PlainTextCell $ws.Cells.Item(109, 9) 'static/uploads/S105.png'

# Row 110: S106
Set-PlainTextCell $ws.Cells.Item(110, 1) 'S106'
Set-PlainTextCell $ws.Cells.Item(110, 2) 'S'
Set-NikCell $ws.Cells.Item(110, 3) '3319012009710001'
Set-PlainTextCell $ws.Cells.Item(110, 4) 'ACHMAD SYAFI''I'
Set-PlainTextCell $ws.Cells.Item(110, 5) '0'
Set-PlainTextCell $ws.Cells.Item(110, 6) 'DUKUH GROGOL RT/RW 006/003 DESA BAKALAN KRAPYAK, KEC KALIWUNGU, KAB KUDUS'
$ws.Cells.Item(110, 7).Value = 0
Set-PlainTextCell $ws.Cells.Item(110, 8) ''
Set-PlainTextCell $ws.Cells.Item(110, 9) 'static/uploads/S106.png'

# Row 111: S107
Set-PlainTextCell $ws.Cells.Item(111, 1) 'S107'
Set-PlainTextCell $ws.Cells.Item(111, 2) 'S'
Set-NikCell $ws.Cells.Item(111, 3) '33190823027700001'
Set-PlainTextCell $ws.Cells.Item(111, 4) 'MOCH RIDWAN'
Set-PlainTextCell $ws.Cells.Item(111, 5) '081325416658'
Set-PlainTextCell $ws.Cells.Item(111, 6) 'RT 02/01 DESAA KARANGMALANG, KEC GEBOG, KAB KUDUS'
$ws.Cells.Item(111, 7).Value = 0
Set-PlainTextCell $ws.Cells.Item(111, 8) ''
Set-PlainTextCell $ws.Cells.Item(111, 9) 'static/uploads/S107.png'

# Row 112: S108
Set-PlainTextCell $ws.Cells.Item(112, 1) 'S108'
Set-PlainTextCell $ws.Cells.Item(112, 2) 'S'
Set-NikCell $ws.Cells.Item(112, 3) '3319060305720002'
Set-PlainTextCell $ws.Cells.Item(112, 4) 'SURYANI'
Set-PlainTextCell $ws.Cells.Item(112, 5) '081391387397'
Set-PlainTextCell $ws.Cells.Item(112, 6) 'RT 02/07 DESA TERBAN, KEC JEKULO, KAB KUDUS'
$ws.Cells.Item(112, 7).Value = 0
Set-PlainTextCell $ws.Cells.Item(112, 8) ''
Set-PlainTextCell $ws.Cells.Item(112, 9) 'static/uploads/S108.png'

# Row 113: S109
Set-PlainTextCell $ws.Cells.Item(113, 1) 'S109'
Set-PlainTextCell $ws.Cells.Item(113, 2) 'S'
Set-NikCell $ws.Cells.Item(113, 3) '3319062710840005'
Set-PlainTextCell $ws.Cells.Item(113, 4) 'JOKO SAPUTRO'
Set-PlainTextCell $ws.Cells.Item(113, 5) '081228768461'
Set-PlainTextCell $ws.Cells.Item(113, 6) 'RT 04/03 DESA SIDOMULYO, KEC JEKULO, KAB KUDUS'
$ws.Cells.Item(113, 7).Value = 0
Set-PlainTextCell $ws.Cells.Item(113, 8) ''
Set-PlainTextCell $ws.Cells.Item(113, 9) 'static/uploads/S109.png'

# Row 114: S110
Set-PlainTextCell $ws.Cells.Item(114, 1) 'S110'
Set-PlainTextCell $ws.Cells.Item(114, 2) 'S'
Set-NikCell $ws.Cells.Item(114, 3) '3319032205790001'
Set-PlainTextCell $ws.Cells.Item(114, 4) 'ABDUL MALIK'
Set-PlainTextCell $ws.Cells.Item(114, 5) '082227937092'
Set-PlainTextCell $ws.Cells.Item(114, 6) 'RT 03/06 DESA LORAM WETAN, KEC JATI, KAB KUDUS'
$ws.Cells.Item(114, 7).Value = 0
Set-PlainTextCell $ws.Cells.Item(114, 8) ''
Set-PlainTextCell $ws.Cells.Item(114, 9) 'static/uploads/S110.png'

# Row 115: S111
Set-PlainTextCell $ws.Cells.Item(115, 1) 'S111'
Set-PlainTextCell $ws.Cells.Item(115, 2) 'S'
Set-NikCell $ws.Cells.Item(115, 3) '3319033012640001'
Set-PlainTextCell $ws.Cells.Item(115, 4) 'ABU BADARI'
Set-PlainTextCell $ws.Cells.Item(115, 5) '085290148365'
Set-PlainTextCell $ws.Cells.Item(115, 6) 'RT 02/04 DESA JATI KULON, KEC JATI, KAB KUDUS'
$ws.Cells.Item(115, 7).Value = 0
Set-PlainTextCell $ws.Cells.Item(115, 8) ''
Set-PlainTextCell $ws.Cells.Item(115, 9) 'static/uploads/S111.png'

# Row 116: S112
Set-PlainTextCell $ws.Cells.Item(116, 1) 'S112'
Set-PlainTextCell $ws.Cells.Item(116, 2) 'S'
Set-NikCell $ws.Cells.Item(116, 3) '3319020211010004'
Set-PlainTextCell $ws.Cells.Item(116, 4) 'MUUHAMMAD NI''MAL MAULANA SISWANTO'
Set-PlainTextCell $ws.Cells.Item(116, 5) '085523784468'
Set-PlainTextCell $ws.Cells.Item(116, 6) 'RT 04/02 DESA KALIPUTU, KEC KOTA, KAB KUDUS'
$ws.Cells.Item(116, 7).Value = 0
Set-PlainTextCell $ws.Cells.Item(116, 8) ''
Set-PlainTextCell $ws.Cells.Item(116, 9) 'static/uploads/S112.png'
